# Updated cryptos list — refresh prices / 1h volume percentages, and
# correct the ranking order for ImmutableX / Fetch.AI (rows 36 and 37
# swap places).
#
# Note: several "Price" cells (column D) are numeric-looking text (e.g.
# "556.97"); Excel auto-converts a bare numeric string typed into a cell
# to a real number, so those are entered with a leading apostrophe to
# force them to stay text, exactly as typing them in the Excel UI would
# require. The Style is then reset to "Normal" so the cell doesn't pick
# up the quote-prefix formatting flag Excel would otherwise remember.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

# --- Row 36 / 37 swap: Fetch.AI now ranks above ImmutableX ---------------
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D36" "0.928"
$ws.Range("E36").Value = "  +5.42%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D37" "1.19"
$ws.Range("E37").Value = "  +5.31%  "

# --- Price (column D) / Volume 1h (column E) refresh ----------------------
Set-TextValue "D2" "59.784.33"
$ws.Range("E2").Value = "  +1.46%  "

Set-TextValue "D3" "2.603.13"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("E4").Value = "  +0.07%  "

Set-TextValue "D5" "556.97"
$ws.Range("E5").Value = "  -1.61%  "

Set-TextValue "D6" "141.73"
$ws.Range("E6").Value = "  -1.10%  "

Set-TextValue "D7" "0.998"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  +0.59%  "

Set-TextValue "D9" "2.625.74"
$ws.Range("E9").Value = "  +1.65%  "

$ws.Range("E10").Value = "  +0.22%  "

$ws.Range("E11").Value = "  +1.77%  "

Set-TextValue "D12" "0.161"
$ws.Range("E12").Value = "  +5.69%  "

$ws.Range("E13").Value = "  +9.01%  "

Set-TextValue "D14" "3.069.56"
$ws.Range("E14").Value = "  +1.31%  "

Set-TextValue "D15" "23.36"
$ws.Range("E15").Value = "  +5.05%  "

Set-TextValue "D16" "59.757.98"
$ws.Range("E16").Value = "  +1.29%  "

Set-TextValue "D17" "0.0000138"
$ws.Range("E17").Value = "  +0.88%  "

Set-TextValue "D18" "2.613.27"
$ws.Range("E18").Value = "  +1.32%  "

Set-TextValue "D19" "4.64"
$ws.Range("E19").Value = "  +3.28%  "

Set-TextValue "D20" "342.89"
$ws.Range("E20").Value = "  +2.08%  "

Set-TextValue "D21" "10.67"
$ws.Range("E21").Value = "  +5.05%  "

Set-TextValue "D22" "6.83"
$ws.Range("E22").Value = "  +10.73%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("E24").Value = "  +14.94%  "

Set-TextValue "D25" "62.54"
$ws.Range("E25").Value = "  -2.18%  "

$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("E27").Value = "  -1.11%  "

Set-TextValue "D28" "7.67"
$ws.Range("E28").Value = "  +6.07%  "

$ws.Range("D29").Value = "0.0₃0783"
$ws.Range("E29").Value = "  +0.48%  "

$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("E32").Value = "  +2.28%  "

Set-TextValue "D33" "158.71"

Set-TextValue "D34" "19.47"
$ws.Range("E34").Value = "  +2.64%  "

$ws.Range("E35").Value = "  +3.72%  "

Set-TextValue "D38" "37.78"
$ws.Range("E38").Value = "  +2.78%  "

Set-TextValue "D39" "1.53"
$ws.Range("E39").Value = "  +2.16%  "

Set-TextValue "D40" "0.847"
$ws.Range("E40").Value = "  -2.75%  "

Set-TextValue "D41" "3.73"
$ws.Range("E41").Value = "  +2.74%  "

Set-TextValue "D42" "294.75"
$ws.Range("E42").Value = "  +1.01%  "

Set-TextValue "D43" "140.77"
$ws.Range("E43").Value = "  +12.88%  "

$ws.Range("E44").Value = "  -0.16%  "

Set-TextValue "D45" "0.0980"
$ws.Range("E45").Value = "  +1.12%  "

$ws.Range("E46").Value = "  +1.23%  "

$ws.Range("E47").Value = "  +3.98%  "

Set-TextValue "D48" "0.0539"
$ws.Range("E48").Value = "  +1.27%  "

$ws.Range("E49").Value = "  +0.08%  "

Set-TextValue "D50" "4.76"
$ws.Range("E50").Value = "  +6.06%  "

$ws.Range("E51").Value = "  +2.89%  "
